# Insert two new weekly-report rows right after row 19 (i.e. at positions
# 20-21), pushing the previously-existing rows 20-46 down to rows 22-48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20:A21").EntireRow.Insert()

# New row 20: Cebollín, Primera, week of 2022-11-29, $/paquete 6 unidades
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44894
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112037
$ws.Range("G20").Value = "Cebollín"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 600
$ws.Range("L20").Value = 700
$ws.Range("M20").Value = 650
$ws.Range("N20").Value = "$/paquete 6 unidades"
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 108
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = "Hortaliza"

# New row 21: Cebollín, Segunda, same week, $/paquete 6 unidades
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44894
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112037
$ws.Range("G21").Value = "Cebollín"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 500
$ws.Range("M21").Value = 500
$ws.Range("N21").Value = "$/paquete 6 unidades"
$ws.Range("O21").Value = "Provincia de Diguillín"
$ws.Range("P21").Value = 83
$ws.Range("Q21").Value = 6
$ws.Range("R21").Value = "Hortaliza"

# Keep the date column's existing number format (YYYY-MM-DD HH:MM:SS) applied
# to the two freshly-inserted date cells, matching the rest of column D.
$ws.Range("D20:D21").NumberFormat = $ws.Range("D22").NumberFormat
